$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of trade data (row 13), matching the formatting of the
# row above it (in particular the date format used in column G).
$ws.Range("G12").Copy() | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("A13").Value = 10012.86
$ws.Range("B13").Value = 10064.19
$ws.Range("C13").Value = 17.8
$ws.Range("D13").Value = 17.89
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = 0.51
$ws.Range("G13").Value = 42620.766319444447
$ws.Range("H13").Value = $false
